$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4481.5386
$ws.Range("I64").Value = 5725
$ws.Range("K64").Value = 5725
$ws.Range("M64").Value = -5477
$ws.Range("H67").Value = 4481.5386
$ws.Range("I67").Value = 5725
$ws.Range("K67").Value = 5725
$ws.Range("M67").Value = -4867
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
$ws.Range("H132").Value = 2428.121
$ws.Range("I132").Value = 2310.5806
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 6931.7418
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -4401.7418
$ws.Range("N132").Value = -17810
$ws.Range("H135").Value = 2073.7368
$ws.Range("I135").Value = 1407.2307
$ws.Range("J135").Value = 3517.8333
$ws.Range("K135").Value = 12665.0763
$ws.Range("L135").Value = 31660.4997
$ws.Range("M135").Value = -10130.0763
$ws.Range("N135").Value = -36730.4997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35970.38
$ws.Range("I2").Value = 49048.953
$ws.Range("K2").Value = 49048.953
$ws.Range("M2").Value = -48935.953
$ws.Range("H107").Value = 36850
$ws.Range("J107").Value = 36850
$ws.Range("L107").Value = 36850
$ws.Range("N107").Value = -44530
$ws.Range("H109").Value = 24500
$ws.Range("J109").Value = 24500
$ws.Range("L109").Value = 24500
$ws.Range("N109").Value = -27274
$ws.Range("H116").Value = 35970.38
$ws.Range("I116").Value = 49048.953
$ws.Range("K116").Value = 49048.953
$ws.Range("M116").Value = -46754.953
$ws.Range("H132").Value = 1365.6052
$ws.Range("I132").Value = 828.5
$ws.Range("J132").Value = 3379.75
$ws.Range("K132").Value = 2485.5
$ws.Range("L132").Value = 10139.25
$ws.Range("M132").Value = 44.5
$ws.Range("N132").Value = -15199.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35970.38
$ws.Range("I3").Value = 49048.953
$ws.Range("K3").Value = 49048.953
$ws.Range("M3").Value = -48934.953
$ws.Range("H22").Value = 321.7647
$ws.Range("I22").Value = 321.7647
$ws.Range("K22").Value = 321.7647
$ws.Range("M22").Value = -148.7647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5009.72
$ws.Range("I58").Value = 933.7857
$ws.Range("J58").Value = 10197.272
$ws.Range("K58").Value = 933.7857
$ws.Range("L58").Value = 10197.272
$ws.Range("M58").Value = -730.7857
$ws.Range("N58").Value = -10603.272
$ws.Range("H62").Value = 2891.625
$ws.Range("I62").Value = 2648.8333
$ws.Range("J62").Value = 3620
$ws.Range("K62").Value = 2648.8333
$ws.Range("L62").Value = 3620
$ws.Range("M62").Value = -2024.8333
$ws.Range("N62").Value = -4868
$ws.Range("H65").Value = 2891.625
$ws.Range("I65").Value = 2648.8333
$ws.Range("J65").Value = 3620
$ws.Range("K65").Value = 13244.1665
$ws.Range("L65").Value = 18100
$ws.Range("M65").Value = -10124.1665
$ws.Range("N65").Value = -24340
$ws.Range("H74").Value = 13483.5
$ws.Range("J74").Value = 13483.5
$ws.Range("L74").Value = 13483.5
$ws.Range("N74").Value = -15231.5
$ws.Range("H77").Value = 13483.5
$ws.Range("J77").Value = 13483.5
$ws.Range("L77").Value = 40450.5
$ws.Range("N77").Value = -49186.5
$ws.Range("H88").Value = 16666.666
$ws.Range("I88").Value = 10000
$ws.Range("K88").Value = 10000
$ws.Range("M88").Value = -9594
$ws.Range("H91").Value = 16666.666
$ws.Range("I91").Value = 10000
$ws.Range("K91").Value = 10000
$ws.Range("M91").Value = -8596
$ws.Range("H92").Value = 11998.5
$ws.Range("J92").Value = 11998.5
$ws.Range("L92").Value = 11998.5
$ws.Range("N92").Value = -16990.5
$ws.Range("H107").Value = 381.25
$ws.Range("I107").Value = 226.66667
$ws.Range("J107").Value = 403.33334
$ws.Range("K107").Value = 226.66667
$ws.Range("L107").Value = 403.33334
$ws.Range("M107").Value = 1693.33333
$ws.Range("N107").Value = -4243.33334
$ws.Range("H122").Value = 1239.875
$ws.Range("I122").Value = 1136.6552
$ws.Range("J122").Value = 1512
$ws.Range("K122").Value = 3409.9656
$ws.Range("L122").Value = 4536
$ws.Range("M122").Value = -959.9655999999995
$ws.Range("N122").Value = -9436
$ws.Range("H132").Value = 1493.3922
$ws.Range("I132").Value = 793.75
$ws.Range("K132").Value = 2381.25
$ws.Range("M132").Value = 148.75
$ws.Range("H134").Value = 1205.9445
$ws.Range("I134").Value = 991.7
$ws.Range("J134").Value = 1473.75
$ws.Range("K134").Value = 2975.1
$ws.Range("L134").Value = 4421.25
$ws.Range("M134").Value = -440.1000000000004
$ws.Range("N134").Value = -9491.25
$ws.Range("H136").Value = 5009.72
$ws.Range("I136").Value = 933.7857
$ws.Range("J136").Value = 10197.272
$ws.Range("K136").Value = 2801.3571
$ws.Range("L136").Value = 30591.816
$ws.Range("M136").Value = -251.3571000000002
$ws.Range("N136").Value = -35691.81600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2125.0588
$ws.Range("I129").Value = 1017
$ws.Range("J129").Value = 3708
$ws.Range("K129").Value = 3051
$ws.Range("L129").Value = 11124
$ws.Range("M129").Value = 1949
$ws.Range("N129").Value = -21124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4998.5
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
$ws.Range("H73").Value = 4998.5
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
$ws.Range("H132").Value = 2549.5
$ws.Range("I132").Value = 2204.2856
$ws.Range("J132").Value = 4966
$ws.Range("K132").Value = 6612.8568
$ws.Range("L132").Value = 14898
$ws.Range("M132").Value = -4082.8568
$ws.Range("N132").Value = -19958

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4832840
$ws.Range("I7").Value = 2308.6667
$ws.Range("J7").Value = 10102511
$ws.Range("K7").Value = 2308.6667
$ws.Range("L7").Value = 10102511
$ws.Range("M7").Value = -2196.6667
$ws.Range("N7").Value = -10102735
$ws.Range("H61").Value = 5907.136
$ws.Range("I61").Value = 7753.5625
$ws.Range("J61").Value = 983.3333
$ws.Range("K61").Value = 7753.5625
$ws.Range("L61").Value = 983.3333
$ws.Range("M61").Value = -7551.5625
$ws.Range("N61").Value = -1387.3333
$ws.Range("H113").Value = 5907.136
$ws.Range("I113").Value = 7753.5625
$ws.Range("J113").Value = 983.3333
$ws.Range("K113").Value = 7753.5625
$ws.Range("L113").Value = 983.3333
$ws.Range("M113").Value = -5583.5625
$ws.Range("N113").Value = -5323.3333
$ws.Range("H126").Value = 4832840
$ws.Range("I126").Value = 2308.6667
$ws.Range("J126").Value = 10102511
$ws.Range("K126").Value = 6926.000100000001
$ws.Range("L126").Value = 30307533
$ws.Range("M126").Value = -4456.000100000001
$ws.Range("N126").Value = -30312473

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1521.5264
$ws.Range("I122").Value = 1529.2142
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4587.642599999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2137.642599999999
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 1198.15
$ws.Range("I132").Value = 915.3333
$ws.Range("J132").Value = 1543.8148
$ws.Range("K132").Value = 2745.9999
$ws.Range("L132").Value = 4631.4444
$ws.Range("M132").Value = -215.9998999999998
$ws.Range("N132").Value = -9691.4444
